# Add a "type" column (char/numeric) to the Names sheet, and make the
# Names sheet the active/selected tab (it was the "Input" sheet before).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Names")

# Header for new column H.
$ws.Range("H1").Value = "type"

# Row 2 (Study_ID / "Study Name") is the character/string field.
$ws.Range("H2").Value = "char"

# Rows 3-33 (group_ID, prepost, Mean, SD, SE, N1, N2, Median, min, max,
# CI%, ulci, llci, q1, q3, pval, patient_data, ...) are numeric fields.
for ($r = 3; $r -le 33; $r++) {
    $ws.Cells.Item($r, 8).Value = "numeric"
}

# Bring the Names sheet to the foreground (moves tabSelected from Input
# to Names and sets activeTab on the workbook view), matching the final
# selection/scroll position recorded in the sheet.
$ws.Activate() | Out-Null

$win = $wb.Windows.Item(1)
$win.ScrollRow = 19
$win.ScrollColumn = 1

$ws.Range("E27").Select() | Out-Null
